# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The error table (rows 2-11, columns B:F) needs to be shifted down by one
# row (the newest period's row of errors is inserted at the top, row 2,
# and every subsequent row takes on the values that used to belong to the
# row above it). The rank/count column G is incremented by one for every
# row, since the sample count (N) used to compute each accumulated error
# grew by one observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for columns B:F, rows 2 through 10,
# before we start overwriting cells, since rows 3:11 will take on the
# values currently held by rows 2:10. (Use .Value2 to read back a plain
# numeric scalar.)
$shiftedValues = @{}
for ($r = 2; $r -le 10; $r++) {
    $rowVals = @()
    for ($c = 2; $c -le 6; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $shiftedValues[$r] = $rowVals
}

# Push rows 2:10 down into rows 3:11 (columns B:F).
for ($r = 10; $r -ge 2; $r--) {
    $vals = $shiftedValues[$r]
    for ($c = 2; $c -le 6; $c++) {
        $ws.Cells.Item($r + 1, $c).Value = $vals[$c - 2]
    }
}

# Write the new, newest-period error values into row 2 (columns B:F).
$newRow2 = @(0.01140061561852912, 0.1285420665309999, 0.03014188771060794, 0.1736141921347674, 0.1779866349240244)
for ($c = 2; $c -le 6; $c++) {
    $ws.Cells.Item(2, $c).Value = $newRow2[$c - 2]
}

# Increment the sample-count column (G) by one for every data row.
for ($r = 2; $r -le 11; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    $ws.Cells.Item($r, 7).Value = $g + 1
}
